# Updated cryptos list data (Price / Volume(1h) columns, plus a
# FraxShare/Aptos row swap) to match the latest coinranking.com snapshot.
#
# NOTE: Price values are stored as text in this sheet (e.g. "1.0000",
# "27.253.72"). Assigning a numeric-looking string straight to .Value
# would let Excel auto-coerce it to a real number (dropping trailing
# zeros, mangling thousand-grouped "27.253.72"-style strings, or even
# flipping tiny values into scientific notation). Prefixing the string
# with a literal leading apostrophe ('' inside a single-quoted
# PowerShell string escapes to one literal ' character) forces Excel to
# keep it as plain text, same as a user typing '1.0000 into the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.253.72'
$ws.Range("E2").Value = '  +1.38%  '

$ws.Range("D3").Value = '''1.909.48'
$ws.Range("E3").Value = '  +2.05%  '

$ws.Range("D4").Value = '''1.0000'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''308.31'
$ws.Range("E5").Value = '  +1.13%  '

$ws.Range("D6").Value = '''1.000'
$ws.Range("E6").Value = '  +0.04%  '

$ws.Range("D7").Value = '''0.5248'
$ws.Range("E7").Value = '  +3.18%  '

$ws.Range("D8").Value = '''0.3790'
$ws.Range("E8").Value = '  +3.50%  '

$ws.Range("D9").Value = '''0.07287'
$ws.Range("E9").Value = '  +1.37%  '

$ws.Range("D10").Value = '''21.33'
$ws.Range("E10").Value = '  +3.17%  '

$ws.Range("D11").Value = '''0.9009'
$ws.Range("E11").Value = '  +0.90%  '

$ws.Range("D12").Value = '''0.07682'
$ws.Range("E12").Value = '  +2.14%  '

$ws.Range("D13").Value = '''1.911.36'
$ws.Range("E13").Value = '  +1.80%  '

$ws.Range("D14").Value = '''95.14'
$ws.Range("E14").Value = '  +0.17%  '

$ws.Range("D15").Value = '''5.278'
$ws.Range("E15").Value = '  +1.01%  '

$ws.Range("E16").Value = '  -0.01%  '

$ws.Range("D17").Value = '''0.000008697'
$ws.Range("E17").Value = '  +2.32%  '

$ws.Range("D18").Value = '''14.54'
$ws.Range("E18").Value = '  +2.47%  '

$ws.Range("D19").Value = '''0.9997'
$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("D20").Value = '''27.314.70'
$ws.Range("E20").Value = '  +1.39%  '

$ws.Range("D21").Value = '''5.089'
$ws.Range("E21").Value = '  +1.52%  '

$ws.Range("D22").Value = '''2.143.95'
$ws.Range("E22").Value = '  +1.63%  '

$ws.Range("D24").Value = '''6.448'
$ws.Range("E24").Value = '  +0.87%  '

$ws.Range("D25").Value = '''2.347'
$ws.Range("E25").Value = '  +12.27%  '

$ws.Range("D26").Value = '''145.93'
$ws.Range("E26").Value = '  -1.55%  '

$ws.Range("D27").Value = '''18.19'
$ws.Range("E27").Value = '  +1.78%  '

$ws.Range("D28").Value = '''1.737'
$ws.Range("E28").Value = '  -2.39%  '

$ws.Range("D29").Value = '''114.98'
$ws.Range("E29").Value = '  +1.53%  '

$ws.Range("D30").Value = '''4.963'
$ws.Range("E30").Value = '  +4.87%  '

$ws.Range("D31").Value = '''4.818'
$ws.Range("E31").Value = '  +2.50%  '

$ws.Range("D32").Value = '''0.09240'
$ws.Range("E32").Value = '  +1.13%  '

$ws.Range("D33").Value = '''0.05077'
$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("D34").Value = '''1.246'
$ws.Range("E34").Value = '  +7.81%  '

$ws.Range("D35").Value = '''0.7940'
$ws.Range("E35").Value = '  +5.85%  '

$ws.Range("D36").Value = '''3.003'
$ws.Range("E36").Value = '  +1.24%  '

$ws.Range("D37").Value = '''3.307'
$ws.Range("E37").Value = '  +2.32%  '

$ws.Range("D38").Value = '''2.613'
$ws.Range("E38").Value = '  +3.32%  '

$ws.Range("D39").Value = '''0.5708'
$ws.Range("E39").Value = '  +2.20%  '

$ws.Range("D40").Value = '''0.02000'
$ws.Range("E40").Value = '  +0.26%  '

$ws.Range("E41").Value = '  +0.16%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").Value = '''6.671'
$ws.Range("E42").Value = '  +0.68%  '

$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '''9.021'
$ws.Range("E43").Value = '  +5.04%  '

$ws.Range("D44").Value = '''119.17'
$ws.Range("E44").Value = '  +2.81%  '

$ws.Range("D45").Value = '''0.1523'
$ws.Range("E45").Value = '  +3.16%  '

$ws.Range("D46").Value = '''0.4870'
$ws.Range("E46").Value = '  +2.35%  '

$ws.Range("D47").Value = '''10.23'
$ws.Range("E47").Value = '  +1.32%  '

$ws.Range("D48").Value = '''1.000'
$ws.Range("E48").Value = '  +0.08%  '

$ws.Range("D49").Value = '''1.613'
$ws.Range("E49").Value = '  +2.82%  '

$ws.Range("D50").Value = '''37.60'
$ws.Range("E50").Value = '  +1.75%  '

$ws.Range("D51").Value = '''64.20'
$ws.Range("E51").Value = '  +1.72%  '
